$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19, shifting existing rows 19:24 down to 20:25
$ws.Rows.Item(19).Insert()

# Populate the new row 19 with the new weekly record
$ws.Range("A19").Value = 7
$ws.Range("B19").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C19").Value = "Ñuble"
$ws.Range("D19").Value = 45176
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 100112039
$ws.Range("G19").Value = "Ciboulette"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 2500
$ws.Range("L19").Value = 2500
$ws.Range("M19").Value = 2500
$ws.Range("N19").Value = "$/docena de atados"
$ws.Range("O19").Value = "Región Metropolitana"
$ws.Range("P19").Value = 833
$ws.Range("Q19").Value = 3
$ws.Range("R19").Value = "Hortaliza"

# Match the style (date number format) used by the other D-column cells
$ws.Range("D19").NumberFormat = $ws.Range("D20").NumberFormat
